# Replace the "select_one yes_no" survey row with an explicit "boolean" type,
# and drop the now-unused "yes_no" choice list from the choices sheet.

$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# --- survey sheet -----------------------------------------------------
# Insert a new row above the existing "text" row (row 6) and fill it in
# with the new boolean question, pushing everything below down by one.
$survey.Rows.Item(6).Insert()
$survey.Range("A6").Value = "boolean"
$survey.Range("B6").Value = "boolean"
$survey.Range("C6").Value = "Boolean"

# The old "select_one yes_no" row (originally row 14) is now at row 15;
# remove it entirely since booleans no longer need a choice list.
$survey.Rows.Item(15).Delete()

# --- choices sheet ------------------------------------------------------
# Remove the "yes_no" choice list rows (3 and 4), plus the blank separator
# row (5) that used to sit between it and the "mealtime" list.
$choices.Rows("3:5").Delete()
$choices.Range("B14").Select()

# Reselect on the survey sheet last so it stays the active tab (as in the
# original workbook) while also landing the cursor where the edit expects.
$survey.Range("F11").Select()
